# Apply payroll recalculation + new employee row per commit "commit theo y ni len"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-9 (C, D, E, F, H columns) ---

# Row 2 - Nguyen Van An
$ws.Range("C2").Value = 2398.8462
$ws.Range("D2").Value = 1500.0
$ws.Range("E2").Value = 200.0
$ws.Range("F2").Value = 3898.8462
$ws.Range("H2").Value = 2198.8462

# Row 3 - Tran Thi Binh
$ws.Range("C3").Value = 2398.8462
$ws.Range("D3").Value = 1500.0
$ws.Range("E3").Value = 0.0
$ws.Range("F3").Value = 4398.846
$ws.Range("H3").Value = 3898.8462

# Row 4 - Le Van Cuong
$ws.Range("C4").Value = 2398.8462
$ws.Range("D4").Value = 0.0
$ws.Range("E4").Value = 0.0
$ws.Range("F4").Value = 2398.8462
$ws.Range("H4").Value = 2398.8462

# Row 5 - Pham Thi Duyen
$ws.Range("C5").Value = 2429.9998
$ws.Range("D5").Value = 1500.0
$ws.Range("E5").Value = 0.0
$ws.Range("F5").Value = 3929.9998
$ws.Range("H5").Value = 4430.0

# Row 6 - Hoang Van Em
$ws.Range("C6").Value = 2429.9998
$ws.Range("D6").Value = 2000.0
$ws.Range("E6").Value = 0.0
$ws.Range("F6").Value = 2429.9998
$ws.Range("H6").Value = 2429.9998

# Row 7 - Ngo Thi Hoa
$ws.Range("C7").Value = 4159.0386
$ws.Range("D7").Value = 2000.0
$ws.Range("E7").Value = 0.0
$ws.Range("F7").Value = 5659.0386
$ws.Range("H7").Value = 5659.0386

# Row 8 - Do Van Khai
$ws.Range("C8").Value = 4159.0386
$ws.Range("D8").Value = 2000.0
$ws.Range("E8").Value = 0.0
$ws.Range("F8").Value = 6159.0386
$ws.Range("H8").Value = 5659.0386

# Row 9 - Bui Thi Lan
$ws.Range("C9").Value = 4159.0386
$ws.Range("D9").Value = 1500.0
$ws.Range("E9").Value = 0.0
$ws.Range("F9").Value = 6159.0386
$ws.Range("H9").Value = 6159.0386

# --- Add new row 10 for a new employee "A" ---
$ws.Range("A10").Value = 9.0
$ws.Range("B10").Value = "A"
$ws.Range("C10").Value = 2429.9998
$ws.Range("D10").Value = 1500.0
$ws.Range("E10").Value = 0.0
$ws.Range("F10").Value = 4430.0
$ws.Range("G10").Value = 0.0
$ws.Range("H10").Value = 2429.9998
